$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("DATA")
$tachesSheet = $wb.Worksheets.Item("Taches")

# Update member names on the DATA sheet
$dataSheet.Range("D5").Value = "Dorian"
$dataSheet.Range("D6").Value = "Nils"
$dataSheet.Range("D7").Value = "Théo"
$dataSheet.Range("D8").ClearContents()

# Update the assignments on the Taches sheet to reflect the renamed members
$tachesSheet.Range("F4").Value = "Nils"
$tachesSheet.Range("F5").Value = "Dorian"

# Restore selections to match the saved state
$tachesSheet.Range("F4").Select()
$dataSheet.Range("D13").Select()
